$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-06-17 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-18 Sunday", 2) | Out-Null

# New values for the 100 table cells, in row-major order (5 columns x 20 rows)
$newValues = @(
    "100×55=",
    "25×71=",
    "73×65=",
    "44×87=",
    "54×98=",
    "48×95=",
    "39×38=",
    "83×13=",
    "80×23=",
    "83×65=",
    "72×79=",
    "98×48=",
    "86×78=",
    "76×74=",
    "55×36=",
    "66×81=",
    "47×71=",
    "98×98=",
    "71×61=",
    "85×95=",
    "73×20=",
    "23×18=",
    "26×34=",
    "18×60=",
    "82×53=",
    "36×13=",
    "50×42=",
    "46×32=",
    "14×42=",
    "51×95=",
    "51×90=",
    "70×67=",
    "50×49=",
    "98×39=",
    "98×91=",
    "74×63=",
    "63×94=",
    "38×48=",
    "29×38=",
    "47×82=",
    "55×12=",
    "52×13=",
    "65×79=",
    "90×49=",
    "28×40=",
    "95×42=",
    "90×17=",
    "40×93=",
    "81×23=",
    "96×11=",
    "51×70=",
    "87×23=",
    "52×73=",
    "23×90=",
    "62×73=",
    "59×57=",
    "25×70=",
    "98×96=",
    "31×39=",
    "56×65=",
    "24×48=",
    "62×100=",
    "62×25=",
    "70×59=",
    "59×62=",
    "100×45=",
    "78×87=",
    "23×80=",
    "81×93=",
    "64×13=",
    "76×56=",
    "51×14=",
    "92×29=",
    "32×76=",
    "69×78=",
    "64×29=",
    "66×91=",
    "18×76=",
    "90×90=",
    "73×35=",
    "99×65=",
    "90×54=",
    "77×93=",
    "59×68=",
    "41×55=",
    "81×41=",
    "39×88=",
    "44×52=",
    "15×56=",
    "28×94=",
    "20×68=",
    "47×56=",
    "48×46=",
    "91×18=",
    "42×21=",
    "17×35=",
    "60×26=",
    "22×99=",
    "66×57=",
    "17×98="
)

$t = $d.Tables.Item(1)
$numCols = $t.Columns.Count
$numRows = $t.Rows.Count
$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Done: updated $idx cells"
